$d = $word.ActiveDocument

# 1) Fix the spelling error "Extracirrcular" -> "Extracurricular"
$d.Content.Find.Execute("Extracirrcular", $true, $false, $false, $false, $false, $true, 1, $false, "Extracurricular", 2) | Out-Null

# 2) Move the "_GoBack" bookmark from after "Education" to inside the GPA
#    run, splitting " GPA in Major (Computer Science)" between "Com" and
#    "puter Science)". Bookmarks.Add with an existing name relocates it,
#    so the stale bookmark after "Education" is removed automatically.
$full = $d.Content
$found = $full.Find.Execute("GPA in Major (Computer Science)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $splitPos = $full.Start + 17   # ... (Com | puter Science)
    $bmRange = $d.Range($splitPos, $splitPos)
    $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
}
